# Rename speaker labels in column D ("Speaker") on Sheet1 to shortened codes.
# Mapping:
#   HILLARY LEWIS-WOLFSEN -> T
#   CAROLYN DOBSON        -> T2
#   Student A             -> S
#   Student B             -> SN
# All other speaker names (e.g. ALYSSA, ASHANK, ANDREW) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$map = @{
    "HILLARY LEWIS-WOLFSEN" = "T"
    "CAROLYN DOBSON"        = "T2"
    "Student A"             = "S"
    "Student B"             = "SN"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
